$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "related work" list in column S (and matching column T entry on row 18)
# gains a new entry ("NeurIPS" / "https://neurips.cc/") inserted at row 21.
# Everything that used to live in S21:T36 shifts down by one row, while the
# B-column headings/labels stay put (only S/T move).
#
# Capture the current (pre-shift) values & bold-formatting of every S/T cell
# from row 21 downward before overwriting anything.
$srcRows = 21..36
$values = @{}
$bolds = @{}
foreach ($r in $srcRows) {
    $values["S$r"] = $ws.Range("S$r").Value()
    $bolds["S$r"] = $ws.Range("S$r").Font.Bold
    $values["T$r"] = $ws.Range("T$r").Value()
    $bolds["T$r"] = $ws.Range("T$r").Font.Bold
}

# Write the shifted-down values (row r -> row r+1), walking from the bottom
# up so we never clobber a value before it has been read (values/bolds were
# already captured above, but keep the safe order regardless).
for ($r = 36; $r -ge 21; $r--) {
    $destS = $r + 1
    $destT = $r + 1

    $sVal = $values["S$r"]
    if ($sVal -eq $null -or $sVal -eq "") {
        $ws.Range("S$destS").Clear()
    } else {
        $ws.Range("S$destS").Value = $sVal
        $ws.Range("S$destS").Font.Bold = $bolds["S$r"]
    }

    $tVal = $values["T$r"]
    if ($tVal -eq $null -or $tVal -eq "") {
        $ws.Range("T$destT").Clear()
    } else {
        $ws.Range("T$destT").Value = $tVal
        $ws.Range("T$destT").Font.Bold = $bolds["T$r"]
    }
}

# Now place the new "NeurIPS" entry on the freshly vacated row 21.
$ws.Range("S21").Value = "NeurIPS"
$ws.Range("S21").Font.Bold = $false
$ws.Range("T21").Value = "https://neurips.cc/"
$ws.Range("T21").Font.Bold = $false

# Update the active-cell selection to match the authored state.
$null = $ws.Range("T21").Select()
